$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Card20")

# --- Row 16 -----------------------------------------------------------
# These columns previously held blank placeholder cells; they now hold the
# literal text "nan" (the sheet's existing convention for "not applicable"),
# matching every other data row.
$ws.Range("B16").Value = "nan"
$ws.Range("C16").Value = "nan"
$ws.Range("D16").Value = "nan"
$ws.Range("E16").Value = "nan"
$ws.Range("F16").Value = "nan"
$ws.Range("G16").Value = "nan"
$ws.Range("H16").Value = "nan"
$ws.Range("I16").Value = "nan"
$ws.Range("J16").Value = "nan"
$ws.Range("K16").Value = "nan"
$ws.Range("M16").Value = "nan"

# --- Row 17 (new service event for Card20) -----------------------------
# Column A repeats the card number, stored as text (matching the rest of
# the "card" column).
$ws.Range("A17").Value = "'20"
$ws.Range("A17").Style = "Normal"

# Columns B-K and M carry no data for this event; write them as blank
# text cells (the sheet's placeholder convention for unfilled columns)
# rather than leaving them completely absent.
$ws.Range("B17").Value = "'"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'"
$ws.Range("E17").Style = "Normal"
$ws.Range("F17").Value = "'"
$ws.Range("F17").Style = "Normal"
$ws.Range("G17").Value = "'"
$ws.Range("G17").Style = "Normal"
$ws.Range("H17").Value = "'"
$ws.Range("H17").Style = "Normal"
$ws.Range("I17").Value = "'"
$ws.Range("I17").Style = "Normal"
$ws.Range("J17").Value = "'"
$ws.Range("J17").Style = "Normal"
$ws.Range("K17").Value = "'"
$ws.Range("K17").Style = "Normal"
$ws.Range("M17").Value = "'"
$ws.Range("M17").Style = "Normal"

# Date, correction note and technician for the new half-yearly maintenance
# event.
$ws.Range("L17").Value = "14\10\2024"
$ws.Range("N17").Value = "تم عمل صيانه نصف سنويه"
$ws.Range("O17").Value = "تيم العمل"
